# B6-PowerPoint.pptx edit
#
# 1) Three tables (on slides 14, 15, 16) get their table style switched
#    from {F7DD2100-5732-47A2-ABDB-4170DB7F37F0} to
#    {1D2AEAAB-E60E-42B0-8F8E-46EBEB479384}.
# 2) The deck's theme palette is switched from the "Red Violet"/"Integral"
#    scheme back to the stock "Office" scheme (dk1/lt1/dk2/lt2/accent1-6/
#    hlink/folHlink), which is what the underlying theme part's colour
#    scheme should carry after the edit.

$p = $ppt.ActivePresentation

# --- 1) Re-style the three tables ------------------------------------------------
$newTableStyle = "{1D2AEAAB-E60E-42B0-8F8E-46EBEB479384}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    $shape = $slide.Shapes.Item(1)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle($newTableStyle)
    }
}

# --- 2) Swap the theme colours back to the stock "Office" palette --------------
# RGB values below are standard OLE COLORREF (0x00BBGGRR) encodings of the
# Office theme hex colours: 000000, FFFFFF, 44546A, E7E6E6, 5B9BD5, ED7D31,
# A5A5A5, FFC000, 4472C4, 70AD47, 0563C1, 954F72 (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink, in that order).
$officeColors = @(
    0,        # dk1       000000
    16777215, # lt1       FFFFFF
    6968388,  # dk2       44546A
    15132391, # lt2       E7E6E6
    13998939, # accent1   5B9BD5
    3243501,  # accent2   ED7D31
    10855845, # accent3   A5A5A5
    49407,    # accent4   FFC000
    12874308, # accent5   4472C4
    4697456,  # accent6   70AD47
    12673797, # hlink     0563C1
    7491477   # folHlink  954F72
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Item($i).RGB = $officeColors[$i - 1]
}
